$wb = $excel.ActiveWorkbook

# --- Resumen sheet: update Maximo (C2) ---
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Range("C2").Value = 693.2600010811979

# --- Solucion sheet: new randomized assignment of Pedido -> Salida ---
$wsSolucion = $wb.Worksheets.Item("Solucion")
$pedidoCol = @("Pedido_26","Pedido_71","Pedido_47","Pedido_8","Pedido_24","Pedido_15","Pedido_21","Pedido_38","Pedido_44","Pedido_7","Pedido_74","Pedido_29","Pedido_78","Pedido_35","Pedido_19","Pedido_52","Pedido_18","Pedido_54","Pedido_63","Pedido_56","Pedido_51","Pedido_28","Pedido_59","Pedido_75","Pedido_57","Pedido_25","Pedido_27","Pedido_33","Pedido_80","Pedido_40","Pedido_73","Pedido_66","Pedido_43","Pedido_60","Pedido_22","Pedido_79","Pedido_36","Pedido_6","Pedido_77","Pedido_23","Pedido_42","Pedido_65","Pedido_62","Pedido_61","Pedido_20","Pedido_37","Pedido_70","Pedido_31","Pedido_30","Pedido_76","Pedido_50","Pedido_3","Pedido_4","Pedido_39","Pedido_12","Pedido_32","Pedido_34","Pedido_68","Pedido_13","Pedido_64","Pedido_1","Pedido_5","Pedido_58","Pedido_55","Pedido_14","Pedido_53","Pedido_41","Pedido_49","Pedido_16","Pedido_2","Pedido_46","Pedido_67","Pedido_69","Pedido_10","Pedido_17","Pedido_72","Pedido_9","Pedido_45","Pedido_11","Pedido_48")
$salidaCol = @("S001","S025","S041","S065","S045","S069","S029","S005","S042","S002","S066","S026","S046","S006","S030","S043","S070","S027","S003","S067","S047","S007","S031","S044","S071","S068","S004","S028","S048","S072","S008","S049","S032","S009","S073","S053","S013","S050","S077","S033","S074","S010","S054","S037","S078","S014","S051","S034","S055","S075","S011","S038","S052","S015","S079","S056","S035","S012","S016","S076","S057","S039","S017","S080","S061","S021","S058","S036","S018","S022","S062","S040","S019","S059","S023","S063","S020","S060","S024","S064")
for ($i = 0; $i -lt $pedidoCol.Length; $i++) {
    $row = $i + 2
    $wsSolucion.Cells.Item($row, 1).Value = $pedidoCol[$i]
    $wsSolucion.Cells.Item($row, 2).Value = $salidaCol[$i]
}

# --- Metricas sheet: update Tiempo per Zona ---
$wsMetricas = $wb.Worksheets.Item("Metricas")
$wsMetricas.Range("B2").Value = 693.2600010811979
$wsMetricas.Range("B3").Value = 522.5721699643204
$wsMetricas.Range("B4").Value = 687.8430641150396
$wsMetricas.Range("B5").Value = 442.4714293437128

Write-Output "edit complete"